$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.487.50"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.845.76"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.40"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5199"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3219"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06772"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.63"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7698"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "1.849.33"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.009"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.91"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007939"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "26.528.19"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "2.090.77"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.612"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.420"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.973"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.40"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.176"
$ws.Range("E26").Value = "  -7.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.676"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.99"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.65"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.154"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08725"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.101"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04810"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.127"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.867"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7167"
$ws.Range("E36").Value = "  +3.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.098"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01781"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.186"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4834"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.16"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8956"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.027"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.604"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4162"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05899"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.071"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.89"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8852"
$ws.Range("E51").Value = "  +3.64%  "
